# Update the division problems in the table to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old="152÷3="; new="262÷5="},
    @{old="700÷5="; new="460÷2="},
    @{old="498÷4="; new="766÷5="},
    @{old="119÷9="; new="750÷5="},
    @{old="533÷9="; new="793÷4="},
    @{old="492÷2="; new="857÷6="},
    @{old="270÷5="; new="486÷9="},
    @{old="941÷3="; new="561÷6="},
    @{old="658÷4="; new="803÷8="},
    @{old="423÷8="; new="414÷4="},
    @{old="553÷3="; new="239÷6="},
    @{old="648÷6="; new="293÷2="},
    @{old="104÷2="; new="547÷5="},
    @{old="299÷4="; new="904÷9="},
    @{old="551÷9="; new="117÷7="},
    @{old="861÷5="; new="291÷4="},
    @{old="551÷6="; new="974÷7="},
    @{old="146÷6="; new="199÷7="},
    @{old="495÷8="; new="450÷4="},
    @{old="897÷5="; new="835÷6="},
    @{old="812÷5="; new="470÷4="},
    @{old="233÷7="; new="356÷9="},
    @{old="233÷2="; new="279÷3="},
    @{old="896÷8="; new="615÷9="},
    @{old="771÷2="; new="884÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
